# Burndown chart sprint update
# - Rename tasks to include owners, shuffle the "Sprint 4 Burndown Chart" title string
# - Record day-5 (column I) effort: feature 1 finished (I8=1) and metrics task ticks to 1 (I9=1)
# - Move the active selection to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task / header text updates -------------------------------------------------
# Order matters for how the shared-string table gets rebuilt, so write these in the
# same order the target workbook lists them (B2, C6, C8, C9, C7).
$ws.Range("B2").Value = "Sprint 4 Burndown Chart"
$ws.Range("C6").Value = "Implement feature 1 (Ricardo)"
$ws.Range("C8").Value = "Upload sprints to github (Todo mundo)"
$ws.Range("C9").Value = "Take metrics (João)"
$ws.Range("C7").Value = "Implement feature 2 (James, João, Francisco, Iago)"

# --- Day 5 (column I) effort entries --------------------------------------------
$ws.Range("I8").Value = 1
$ws.Range("I9").Value = 1

# --- Selection mirrors the saved cursor position in the target file -------------
$ws.Range("C9").Select() | Out-Null
